$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Quentin Grimes", "SG,SF", "Philadelphia 76ers"),
    @("Daniel Gafford", "PF,C", "Dallas Mavericks"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Gary Trent Jr.", "PG,SG,SF", "Milwaukee Bucks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Mark Williams", "C", "Charlotte Hornets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
